# Swap the order of "Recorded By" names in column G from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every row where that exact value is present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Output "Rows updated: $changed"
